$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values are written as literal text (not auto-converted to numbers)
# by temporarily forcing a text number format on the whole data range, then clearing
# the formatting afterwards so the resulting style matches the original (no explicit
# style index on these cells).
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '53.879.40'
$ws.Range("E2").Value = '  -11.07%  '

$ws.Range("D3").Value = '2.330.87'
$ws.Range("E3").Value = '  -19.68%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").Value = '443.51'
$ws.Range("E5").Value = '  -15.95%  '

$ws.Range("D6").Value = '125.57'
$ws.Range("E6").Value = '  -12.46%  '

$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.36%  '

$ws.Range("E8").Value = '  -14.31%  '

$ws.Range("D9").Value = '2.330.42'
$ws.Range("E9").Value = '  -19.92%  '

$ws.Range("D10").Value = '5.35'
$ws.Range("E10").Value = '  -11.37%  '

$ws.Range("D11").Value = '0.0917'
$ws.Range("E11").Value = '  -15.04%  '

$ws.Range("D12").Value = '0.308'
$ws.Range("E12").Value = '  -14.89%  '

$ws.Range("E13").Value = '  -3.35%  '

$ws.Range("D14").Value = '2.683.31'
$ws.Range("E14").Value = '  -21.29%  '

$ws.Range("D15").Value = '53.905.94'
$ws.Range("E15").Value = '  -11.03%  '

$ws.Range("D16").Value = '18.81'
$ws.Range("E16").Value = '  -17.41%  '

$ws.Range("E17").Value = '  -14.16%  '

$ws.Range("D18").Value = '2.348.46'
$ws.Range("E18").Value = '  -19.32%  '

$ws.Range("D19").Value = '3.93'
$ws.Range("E19").Value = '  -21.96%  '

$ws.Range("D20").Value = '298.40'
$ws.Range("E20").Value = '  -17.51%  '

$ws.Range("D21").Value = '9.14'
$ws.Range("E21").Value = '  -21.88%  '

$ws.Range("E22").Value = '  -0.17%  '

$ws.Range("D23").Value = '5.59'
$ws.Range("E23").Value = '  -1.67%  '

$ws.Range("D24").Value = '5.38'
$ws.Range("E24").Value = '  -18.97%  '

$ws.Range("D25").Value = '55.62'
$ws.Range("E25").Value = '  -14.27%  '

$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.38%  '

$ws.Range("D27").Value = '0.151'
$ws.Range("E27").Value = '  -16.23%  '

$ws.Range("D28").Value = '0.366'
$ws.Range("E28").Value = '  -19.44%  '

$ws.Range("D29").Value = '6.94'
$ws.Range("E29").Value = '  -11.68%  '

$ws.Range("D30").Value = '0.997'
$ws.Range("E30").Value = '  -0.21%  '

$ws.Range("E31").Value = '  -17.80%  '

$ws.Range("D32").Value = '145.98'
$ws.Range("E32").Value = '  -4.12%  '

$ws.Range("D33").Value = '17.21'
$ws.Range("E33").Value = '  -12.91%  '

$ws.Range("E34").Value = '  -19.63%  '

$ws.Range("D35").Value = '4.65'
$ws.Range("E35").Value = '  -16.74%  '

$ws.Range("D36").Value = '3.53'
$ws.Range("E36").Value = '  -19.38%  '

$ws.Range("D37").Value = '0.831'
$ws.Range("E37").Value = '  -17.49%  '

$ws.Range("E38").Value = '  -16.73%  '

$ws.Range("D39").Value = '33.30'
$ws.Range("E39").Value = '  -11.66%  '

$ws.Range("D40").Value = '0.995'
$ws.Range("E40").Value = '  -0.31%  '

$ws.Range("D41").Value = '10.27'
$ws.Range("E41").Value = '  -0.47%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.939.52'
$ws.Range("E42").Value = '  -15.56%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '3.13'
$ws.Range("E43").Value = '  -16.02%  '

$ws.Range("D44").Value = '1.20'
$ws.Range("E44").Value = '  -18.85%  '

$ws.Range("E45").Value = '  -15.10%  '

$ws.Range("D46").Value = '0.522'
$ws.Range("E46").Value = '  -19.44%  '

$ws.Range("D47").Value = '0.0209'
$ws.Range("E47").Value = '  -11.90%  '

$ws.Range("D48").Value = '0.0830'
$ws.Range("E48").Value = '  -10.13%  '

$ws.Range("D49").Value = '4.03'
$ws.Range("E49").Value = '  -19.38%  '

$ws.Range("D50").Value = '15.78'
$ws.Range("E50").Value = '  -22.97%  '

$ws.Range("E51").Value = '  -3.56%  '

# Strip the temporary text formatting so styling matches the original workbook
$dRange.ClearFormats()
